# Timing issue fix - keywords, updated tc1,2 in ubc01
# The 'CasesTab' query (cell B2) is rewritten to drop the trailing
# `coalesce(co.cohort_description, '') AS `Cohort`` return column
# (and its now-dangling trailing comma on the prior line).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCasesTabQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC01'] and diag.stage_of_disease in [ 'T2N0M0', 'Not Applicable'] OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $newCasesTabQuery

# Row 2 no longer needs as much vertical space once the Cohort column
# line is gone; match rows 3/4's wrapped-text height.
$ws.Rows(2).RowHeight = 290

# Author's view scrolled up / selected B2 instead of B4.
$ws.Range("B2").Select()
